# C5-PowerPoint.pptx edit:
#   The table on the "SOURCES OF FINANCE" slide was tagged with a table
#   style GUID ({3D815F76-...}) that doesn't match any style defined in
#   this deck's tableStyles.xml. Re-point it at the intended built-in
#   PowerPoint table style ({2C9762C4-1503-4E59-B02A-E1500F019639}).

$p = $ppt.ActivePresentation

$oldStyleId = "{3D815F76-348E-4888-AE8F-0BEAEC80FC7C}"
$newStyleId = "{2C9762C4-1503-4E59-B02A-E1500F019639}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
